$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.096649527549744
$ws.Range("B1").Value = 1.950227618217468
$ws.Range("D1").Value = 1.053438425064087
$ws.Range("E1").Value = 1.116754531860352
